$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 535 ("A_TYPEAREA_RECREATIO" ...),
# shifting the existing rows (formerly 535-539) down to 537-541.
$ws.Rows.Item(535).Insert()
$ws.Rows.Item(535).Insert()

# The newly inserted rows don't carry the surrounding data-row formatting,
# so copy it over from the row just below (the original row 535, now at 537).
$ws.Range("A537:D537").Copy()
$ws.Range("A535:D536").PasteSpecial(-4122)

# Fill in the two new dictionary entries (row 535 then row 536).
$ws.Range("A535").Value = "A_TARIF_FLAECHE"
$ws.Range("B535").Value = "K_TARIF"
$ws.Range("C535").Value = "Branchen-/Flächentarifvertrag"
$ws.Range("D535").Value = "XXXBranchen-/Flächentarifvertrag"

$ws.Range("A536").Value = "A_TARIF_HAUS"
$ws.Range("B536").Value = "K_TARIF"
$ws.Range("C536").Value = "Firmen-/Haustarifvertrag"
$ws.Range("D536").Value = "XXXFirmen-/Haustarifvertrag"
